$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("data")

$ws.Range("D5").Value = 2.367899673874997
$ws.Range("D6").Value = 0.07002404621610055
$ws.Range("D7").Value = -0.3559151560458474
$ws.Range("D8").Value = 0.2619431526596601
$ws.Range("D9").Value = 2.498394558263731
$ws.Range("D10").Value = 0.2530015216472968
$ws.Range("D11").Value = 2.42295084177452
$ws.Range("D12").Value = 0.3224670496789456
$ws.Range("D13").Value = 0.4435958093950659
$ws.Range("D14").Value = 0.1936870564379858
$ws.Range("D15").Value = 0.193282093897509
$ws.Range("D16").Value = 0.219220690029325
$ws.Range("D17").Value = -0.1105923956464526
$ws.Range("D18").Value = 0.005863903358051403
$ws.Range("D19").Value = 0.4143630902555517
$ws.Range("D20").Value = 0.5139108704465918
$ws.Range("D21").Value = 0.1719286908784542
$ws.Range("D22").Value = 0.03066205678034255
